# Gate_Planning.xlsx — "Added walking distance to passport control"
#
# Adds a new "Passport Control" worksheet (after "Piers") listing the
# walking distance from each pier to passport control, and updates the
# selections / active sheet left behind by the editing session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "Passport Control" sheet, positioned after "Piers"
# ---------------------------------------------------------------------
$piers = $wb.Worksheets.Item("Piers")
$passport = $wb.Worksheets.Add($null, $piers)
$passport.Name = "Passport Control"

# Copy the header formatting (bold font + medium bottom border) from the
# "Piers" sheet's header row so the new sheet matches the existing look.
$piers.Range("A1:B1").Copy() | Out-Null
$passport.Range("A1:B1").PasteSpecial(-4122) | Out-Null
$passport.Application.CutCopyMode = $false

# Header row
$passport.Range("A1").Value = "Pier"
$passport.Range("B1").Value = "Distance"

# Data rows: walking distance (meters) from each pier to passport control
$passport.Range("A2").Value = "A"
$passport.Range("B2").Value = 75
$passport.Range("A3").Value = "B"
$passport.Range("B3").Value = 75
$passport.Range("A4").Value = "C"
$passport.Range("B4").Value = 75
$passport.Range("A5").Value = "H"
$passport.Range("B5").Value = 150

# ---------------------------------------------------------------------
# 2. Leftover UI state from the editing session: selections per sheet
# ---------------------------------------------------------------------
$flightSchedule = $wb.Worksheets.Item("Flight Schedule")
$flightSchedule.Activate()
$flightSchedule.Range("G50").Select() | Out-Null

$airlines = $wb.Worksheets.Item("Airlines")
$airlines.Activate()
$airlines.Range("A11").Select() | Out-Null

$transfers = $wb.Worksheets.Item("Transfers")
$transfers.Activate()
$transfers.Range("E13").Select() | Out-Null

$gates = $wb.Worksheets.Item("Gates")
$gates.Activate()
$gates.Range("G9").Select() | Out-Null

$passport.Activate()
$passport.Range("J28").Select() | Out-Null

# "Piers" is the sheet that ends up active/selected in the saved workbook
$piers.Activate()
$piers.Range("J29").Select() | Out-Null
